# Apply crypto price/volume update (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B / C / E updates: plain text, safe to assign directly ---
$ws.Range("E2").Value = "  +3.71%  "
$ws.Range("E3").Value = "  +2.73%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("E5").Value = "  +4.21%  "
$ws.Range("E6").Value = "  +6.60%  "
$ws.Range("E7").Value = "  +1.98%  "
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("E9").Value = "  +5.43%  "
$ws.Range("E10").Value = "  +3.93%  "
$ws.Range("E11").Value = "  +1.95%  "
$ws.Range("E12").Value = "  +1.46%  "
$ws.Range("E13").Value = "  +2.94%  "
$ws.Range("E14").Value = "  +3.42%  "
$ws.Range("E15").Value = "  +2.76%  "
$ws.Range("E16").Value = "  +3.61%  "
$ws.Range("E17").Value = "  +5.20%  "
$ws.Range("E18").Value = "  +3.56%  "
$ws.Range("E19").Value = "  +3.62%  "
$ws.Range("E20").Value = "  +2.29%  "
$ws.Range("E21").Value = "  +2.59%  "
$ws.Range("E22").Value = "  +1.55%  "
$ws.Range("E23").Value = "  +2.60%  "
$ws.Range("E24").Value = "  +4.40%  "
$ws.Range("E25").Value = "  +2.58%  "
$ws.Range("E26").Value = "  -0.07%  "
$ws.Range("E27").Value = "  +2.12%  "
$ws.Range("E28").Value = "  -4.28%  "
$ws.Range("E29").Value = "  +4.38%  "
$ws.Range("E30").Value = "  +5.76%  "
$ws.Range("E31").Value = "  +1.23%  "
$ws.Range("B32").Value = "Celestia"
$ws.Range("C32").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("E32").Value = "  +12.66%  "
$ws.Range("B33").Value = "Kaspa"
$ws.Range("C33").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("E33").Value = "  +16.40%  "
$ws.Range("E34").Value = "  +3.50%  "
$ws.Range("E35").Value = "  +0.22%  "
$ws.Range("E36").Value = "  +6.26%  "
$ws.Range("B37").Value = "RenderToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("E37").Value = "  +4.36%  "
$ws.Range("B38").Value = "ARBITRUM"
$ws.Range("C38").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("E38").Value = "  +2.98%  "
$ws.Range("E39").Value = "  +3.85%  "
$ws.Range("E40").Value = "  +4.33%  "
$ws.Range("E41").Value = "  +1.08%  "
$ws.Range("E42").Value = "  +0.63%  "
$ws.Range("E43").Value = "  -6.88%  "
$ws.Range("E44").Value = "  +3.24%  "
$ws.Range("E45").Value = "  +0.73%  "
$ws.Range("E47").Value = "  +8.62%  "
$ws.Range("E48").Value = "  +4.90%  "
$ws.Range("E49").Value = "  +11.14%  "
$ws.Range("E50").Value = "  +3.51%  "
$ws.Range("E51").Value = "  +2.44%  "

# --- Column D updates: force Text format so numeric-looking strings
#     ("19.55", "6.40", etc.) are preserved verbatim instead of being
#     auto-coerced to numbers (which would e.g. drop trailing zeros).
#     Style is reset back to Normal afterwards so no new number format
#     is left applied to the cell.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "44.512.57"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.427.40"
$ws.Range("D3").Style = "Normal"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "314.52"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "101.72"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.514"
$ws.Range("D7").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.513"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.26"
$ws.Range("D10").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.77"
$ws.Range("D13").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.806.43"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.441.78"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.838"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "44.400.16"
$ws.Range("D18").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.40"
$ws.Range("D20").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.90"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "241.04"
$ws.Range("D23").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.19"
$ws.Range("D27").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.64"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "33.33"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "48.47"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "19.55"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.122"
$ws.Range("D33").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0766"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.55"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.90"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.89"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "126.94"
$ws.Range("D40").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "21.88"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.15"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0288"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.948.71"
$ws.Range("D45").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.64"
$ws.Range("D48").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "53.42"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "73.71"
$ws.Range("D51").Style = "Normal"
